$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($r = 2; $r -le 500; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cur = $cell.Value2()
    if ($cur -eq 45202) {
        $cell.Value2 = 45203
    }
}
